$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Next period (release date)" value for "Job adverts by occupation" (row 13)
# The Textkernel job-adverts data source date note changes from a provisional date
# to "TBC*" since the actual date is not yet known.
$ws.Range("D13").Value = "TBC*"

# Update the saved view/selection state: scroll so row 3 is at the top
# and select cell B3 (matches the author's saved cursor position).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B3").Select()

$wb.Save()
